$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell currently stores plain text (inlineStr) that looks numeric
# (e.g. "275.26", "0.61%"). Excel auto-converts a numeric/percent-looking string
# assigned via .Value into a real number, so we force the cell format to Text ("@")
# before writing, keeping the values as literal strings exactly like the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "275.29"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.63%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.24%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.844"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.61%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.931"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.40%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.214"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.81%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8764"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.11%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1516"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.42%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05063"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.51%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07528"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.95%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02962"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.79%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08987"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.54%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001562"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.81%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006434"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.38%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006190"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.83%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.469"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.49%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.31%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.51%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.10%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.909"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.27%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04408"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.08%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.28%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003850"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-12.60%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.05%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "14.01%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04128"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006767"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.85%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.65%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002071"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.85%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005182"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.37%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.695"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-36.42%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.02002"
